# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# Collapses the expanded "dateTime"/body-detail rows on the Body, 200, 204,
# 400 sheets down to a single schema-reference row, and adds the missing
# schema-reference row (errorResponse / errorResponse1) on the 401/403/404/
# 429/500 error sheets.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param(
        $ws,
        [string]$section,
        [string]$name
    )

    $ws.Range("A3").Value = $section
    $ws.Range("B3").Value = $name
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $name
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# --- Body: rows 3-12 (expanded request body) collapse into one schema row ---
$ws = $wb.Worksheets.Item("Body")
$ws.Range("A4:O12").Clear()
Set-SchemaRow $ws "body" "setDefaultAgenda.211207Request"

# --- 200: rows 3-4 (expanded response body) collapse into one schema row ---
$ws = $wb.Worksheets.Item("200")
$ws.Range("A4:O4").Clear()
Set-SchemaRow $ws "content" "setDefaultAgenda.211207Response"

# --- 204: add the missing schema row ---
$ws = $wb.Worksheets.Item("204")
Set-SchemaRow $ws "content" "setDefaultAgenda.211207Response"

# --- 400: rows 3-6 (expanded error body) collapse into one schema row ---
$ws = $wb.Worksheets.Item("400")
$ws.Range("A4:O6").Clear()
Set-SchemaRow $ws "content" "errorResponse"

# --- 401 / 403 / 404 / 429 / 500: add the missing schema row ---
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Set-SchemaRow $ws "content" "errorResponse1"
}
